# notas da prova 3 de ES
# Fill in the "P3 (30)" grades (column D) for every student and refresh the
# related cell formatting (columns D, F, G) to match the new layout, plus the
# widened column H and the updated active-window selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- P3 (column D) grades, row 2 through row 30 -----------------------------
$p3 = @{
    2 = 30;  3 = 30;  4 = 25;  5 = 25;  6 = 25;  7 = 25;  8 = 25;
    9 = 30; 10 = 30; 11 = 30; 12 = 20; 13 = 30; 14 = 25; 15 = 25;
   16 = 30; 17 = 30; 18 = 30; 19 = 25; 20 = 30; 21 = 25; 22 = 30;
   23 = 30; 24 = 25; 25 = 25; 26 = 30; 27 = 25; 28 = 30; 29 = 30; 30 = 30
}

foreach ($row in $p3.Keys) {
    $ws.Cells.Item($row, 4).Value = $p3[$row]
}

# --- Formatting refresh ------------------------------------------------------
# D2:D15 -> Verdana 11, centered (rows for the first batch of grades)
$dTop = $ws.Range("D2:D15")
$dTop.Font.Name = "Verdana"
$dTop.Font.Size = 11
$dTop.HorizontalAlignment = -4108   # xlCenter

# D16:D30 -> Arial 11, centered (second batch)
$dBottom = $ws.Range("D16:D30")
$dBottom.Font.Name = "Arial"
$dBottom.Font.Size = 11
$dBottom.HorizontalAlignment = -4108   # xlCenter

# F2:F30 -> Arial 11 (general alignment, unchanged)
$fCol = $ws.Range("F2:F30")
$fCol.Font.Name = "Arial"
$fCol.Font.Size = 11

# G2:G30 -> Verdana 11, not bold, centered
$gCol = $ws.Range("G2:G30")
$gCol.Font.Name = "Verdana"
$gCol.Font.Size = 11
$gCol.Font.Bold = $false
$gCol.HorizontalAlignment = -4108   # xlCenter

# --- Column H is a little wider now ------------------------------------------
$ws.Columns("H").ColumnWidth = 8.720833333333334

# --- Active window moved down / selection moved to I5 -----------------------
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I5").Select()
